# Error Calculations and Plots
# Apply the edits described by the diff:
#  - Remove two rows that were dropped from the data (old row 26 "RM 232" and
#    old row 28 "SC 92"), which shifts all following rows up and changes the
#    sheet dimension from A1:F35 to A1:F33.
#  - Update a number of individual cells (values that became newly present,
#    or values that were cleared to become "missing" / blank) to match the
#    final state of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that are no longer present in the final data ---
# Deleting row 26 first ("RM 232") shifts "SC 92" up from row 28 to row 27,
# so it is then deleted as the (new) row 27.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Apply the individual cell value changes (post row-shift row numbers) ---
$ws.Range("F2").Value = ""
$ws.Range("F5").Value = 17.66
$ws.Range("E6").Value = -5.7
$ws.Range("F6").Value = 16.43
$ws.Range("E8").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("E12").Value = -5.3
$ws.Range("E14").Value = ""
$ws.Range("E17").Value = -7.3
$ws.Range("E18").Value = -8.5
$ws.Range("E19").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("E23").Value = -7
$ws.Range("F24").Value = 16.78
$ws.Range("B27").Value = -20.4
$ws.Range("E27").Value = ""
$ws.Range("B28").Value = ""
$ws.Range("F28").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("B30").Value = -19.7
$ws.Range("F30").Value = 16.89
$ws.Range("B32").Value = ""
